$wb = $excel.ActiveWorkbook

# The "展览" (exhibitions) sheet and the "全部类型" (all types) sheet contain
# duplicated rows of convention data; both need the same "想去人数" (interest
# count) bumps applied to them.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F7").Value = 1254
    $ws.Range("F17").Value = 300
    $ws.Range("F25").Value = 333
    $ws.Range("F26").Value = 4159
    $ws.Range("F32").Value = 527
    $ws.Range("F36").Value = 136
}
